$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.380.45"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.777.50"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'651.10"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").Value = "'166.04"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "3.777.37"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").Value = "'0.457"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "'6.88"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("E13").Value = "  -4.99%  "
$ws.Range("D14").Value = "'34.99"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").Value = "4.412.90"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").Value = "3.767.06"
$ws.Range("E16").Value = "  -3.98%  "
$ws.Range("D17").Value = "69.293.20"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "'17.81"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'7.02"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "'467.54"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'9.59"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").Value = "'0.709"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -5.47%  "
$ws.Range("D25").Value = "'81.86"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").Value = "'12.37"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").Value = "'10.36"
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").Value = "'2.11"
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "3.926.00"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").Value = "'2.71"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "'2.27"
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").Value = "'7.17"
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("D34").Value = "'28.68"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("D35").Value = "'0.173"
$ws.Range("E35").Value = "  +15.15%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "3.730.23"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'8.84"
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").Value = "'5.84"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").Value = "'3.25"
$ws.Range("E41").Value = "  -6.57%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'0.957"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D45").Value = "'45.23"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").Value = "'1.99"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").Value = "'156.28"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "'47.33"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'0.297"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "'8.37"
$ws.Range("E51").Value = "  -1.04%  "
